$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'58.523.51"
$ws.Cells.Item(3, 4).Value = "'2.626.84"
$ws.Cells.Item(3, 5).Value = "'  +1.05%  "
$ws.Cells.Item(4, 5).Value = "'  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'535.06"
$ws.Cells.Item(5, 5).Value = "'  -0.24%  "
$ws.Cells.Item(6, 4).Value = "'143.09"
$ws.Cells.Item(6, 5).Value = "'  +1.35%  "
$ws.Cells.Item(7, 5).Value = "'  -0.03%  "
$ws.Cells.Item(8, 4).Value = "'0.568"
$ws.Cells.Item(8, 5).Value = "'  +0.35%  "
$ws.Cells.Item(10, 5).Value = "'  -1.72%  "
$ws.Cells.Item(11, 4).Value = "'0.334"
$ws.Cells.Item(11, 5).Value = "'  -0.10%  "
$ws.Cells.Item(12, 5).Value = "'  +0.97%  "
$ws.Cells.Item(13, 4).Value = "'3.093.36"
$ws.Cells.Item(13, 5).Value = "'  +1.18%  "
$ws.Cells.Item(14, 4).Value = "'58.463.44"
$ws.Cells.Item(14, 5).Value = "'  -1.41%  "
$ws.Cells.Item(15, 4).Value = "'20.77"
$ws.Cells.Item(15, 5).Value = "'  +0.43%  "
$ws.Cells.Item(16, 4).Value = "'2.612.17"
$ws.Cells.Item(16, 5).Value = "'  +0.03%  "
$ws.Cells.Item(17, 5).Value = "'  -0.85%  "
$ws.Cells.Item(18, 4).Value = "'4.40"
$ws.Cells.Item(18, 5).Value = "'  +1.08%  "
$ws.Cells.Item(19, 4).Value = "'334.68"
$ws.Cells.Item(19, 5).Value = "'  -1.81%  "
$ws.Cells.Item(20, 4).Value = "'10.14"
$ws.Cells.Item(20, 5).Value = "'  +0.53%  "
$ws.Cells.Item(21, 5).Value = "'  -2.29%  "
$ws.Cells.Item(22, 4).Value = "'0.999"
$ws.Cells.Item(22, 5).Value = "'  -0.07%  "
$ws.Cells.Item(23, 4).Value = "'66.28"
$ws.Cells.Item(23, 5).Value = "'  -1.74%  "
$ws.Cells.Item(24, 4).Value = "'0.416"
$ws.Cells.Item(24, 5).Value = "'  +1.80%  "
$ws.Cells.Item(27, 4).Value = "'7.11"
$ws.Cells.Item(27, 5).Value = "'  -1.47%  "
$ws.Cells.Item(28, 4).Value = "'0.0₃0737"
$ws.Cells.Item(28, 5).Value = "'  -0.72%  "
$ws.Cells.Item(29, 5).Value = "'  -0.01%  "
$ws.Cells.Item(30, 5).Value = "'  -0.97%  "
$ws.Cells.Item(31, 4).Value = "'5.87"
$ws.Cells.Item(31, 5).Value = "'  +1.12%  "
$ws.Cells.Item(32, 4).Value = "'18.79"
$ws.Cells.Item(32, 5).Value = "'  -0.09%  "
$ws.Cells.Item(33, 4).Value = "'150.31"
$ws.Cells.Item(34, 4).Value = "'3.90"
$ws.Cells.Item(34, 5).Value = "'  -1.70%  "
$ws.Cells.Item(35, 4).Value = "'37.20"
$ws.Cells.Item(35, 5).Value = "'  +0.15%  "
$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).Value = "'1.10"
$ws.Cells.Item(36, 5).Value = "'  -0.43%  "
$ws.Cells.Item(37, 2).Value = "SuiNetwork"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Cells.Item(37, 4).Value = "'0.850"
$ws.Cells.Item(37, 5).Value = "'  +2.30%  "
$ws.Cells.Item(38, 5).Value = "'  -3.32%  "
$ws.Cells.Item(39, 4).Value = "'0.810"
$ws.Cells.Item(39, 5).Value = "'  -1.72%  "
$ws.Cells.Item(40, 4).Value = "'3.57"
$ws.Cells.Item(40, 5).Value = "'  +0.96%  "
$ws.Cells.Item(41, 4).Value = "'280.94"
$ws.Cells.Item(41, 5).Value = "'  +3.19%  "
$ws.Cells.Item(42, 5).Value = "'  -0.07%  "
$ws.Cells.Item(43, 5).Value = "'  -0.26%  "
$ws.Cells.Item(44, 4).Value = "'10.68"
$ws.Cells.Item(44, 5).Value = "'  -0.71%  "
$ws.Cells.Item(45, 5).Value = "'  +1.47%  "
$ws.Cells.Item(46, 4).Value = "'19.02"
$ws.Cells.Item(46, 5).Value = "'  +3.07%  "
$ws.Cells.Item(47, 4).Value = "'0.0937"
$ws.Cells.Item(47, 5).Value = "'  -1.66%  "
$ws.Cells.Item(48, 5).Value = "'  +0.46%  "
$ws.Cells.Item(49, 4).Value = "'1.947.70"
$ws.Cells.Item(49, 5).Value = "'  +0.18%  "
$ws.Cells.Item(50, 5).Value = "'  -1.11%  "
$ws.Cells.Item(51, 4).Value = "'17.85"
$ws.Cells.Item(51, 5).Value = "'  -3.97%  "
